$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the party identification text in the opening paragraph
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Emanuel pessoa física, CPF ou CNPJ: 123,   Endereço Residencial ou Comercial: 123, doravante",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Fernandinho pessoa física, CPF ou CNPJ: 55555555555,   Endereço Residencial ou Comercial: Rua Jiboia, doravante",
    2)

# ---------------------------------------------------------------------------
# 2. Horário da Montagem: 07:00 -> 13:27
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "07:00", $true, $false, $false, $false, $false, $true, 1, $false,
    "13:27", 2)

# ---------------------------------------------------------------------------
# 3. Desmontagem: 08:00 -> 14:27
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "08:00", $true, $false, $false, $false, $false, $true, 1, $false,
    "14:27", 2)

# ---------------------------------------------------------------------------
# 4. Paragraph-formatting rework around CLÁUSULA 4 / 4.1
#    - "CLÁUSULA 4 – DAS OBRIGAÇÕES DA CONTRATANTE:" paragraph gains
#      spacing-before (150 twips = 7.5 pt).
#    - "4.1 – O CONTRATANTE compromete-se a: " paragraph loses its
#      spacing-after pPr entirely and loses the bold run formatting.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*CLÁUSULA 4 – DAS OBRIGAÇÕES DA CONTRATANTE*") {
        $p.SpaceBefore = 7.5
    }
}

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*4.1 – O CONTRATANTE compromete-se a*") {
        # Drop the paragraph's own spacing (pPr) by resetting its mark.
        $markRange = $d.Range($p.Range.End - 1, $p.Range.End)
        $markRange.Delete()
        $markRange.InsertParagraphAfter()
    }
}

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*4.1 – O CONTRATANTE compromete-se a*") {
        # Drop the bold run formatting by replacing the run content outright.
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $r.Delete()
        $r.InsertAfter("4.1 – O CONTRATANTE compromete-se a: ")
    }
}

# ---------------------------------------------------------------------------
# 5. Closing date paragraph gains spacing before/after (150 twips = 7.5 pt)
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*São Paulo, 06 de 10 de 2024*") {
        $p.SpaceBefore = 7.5
        $p.SpaceAfter = 7.5
    }
}
